$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 137
$ws.Range("F3").Value = 1323
$ws.Range("F5").Value = 1008
$ws.Range("F6").Value = 1787
$ws.Range("F7").Value = 552
$ws.Range("F8").Value = 1191
$ws.Range("F12").Value = 289
$ws.Range("F15").Value = 683
$ws.Range("F16").Value = 164
$ws.Range("F18").Value = 26
$ws.Range("F21").Value = 140
$ws.Range("F23").Value = 31
$ws.Range("F25").Value = 148
$ws.Range("F27").Value = 870
$ws.Range("F29").Value = 158
$ws.Range("F30").Value = 40
$ws.Range("F33").Value = 14

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 317
$ws.Range("F5").Value = 14

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 308

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 308
$ws.Range("F3").Value = 137
$ws.Range("F4").Value = 1323
$ws.Range("F6").Value = 1008
$ws.Range("F7").Value = 1787
$ws.Range("F9").Value = 1191
$ws.Range("F14").Value = 289
$ws.Range("F17").Value = 683
$ws.Range("F18").Value = 164
$ws.Range("F21").Value = 26
$ws.Range("F22").Value = 317
$ws.Range("F24").Value = 14
$ws.Range("F29").Value = 140
$ws.Range("F31").Value = 31
$ws.Range("F33").Value = 148
$ws.Range("F35").Value = 870
$ws.Range("F39").Value = 158
$ws.Range("F40").Value = 40
$ws.Range("F45").Value = 14
